$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue $ws "D2" '28.040.23'
Set-TextValue $ws "E2" '  -3.93%  '

# Row 3
Set-TextValue $ws "D3" '1.745.05'
Set-TextValue $ws "E3" '  -4.48%  '

# Row 4
Set-TextValue $ws "E4" '  -0.25%  '

# Row 5
Set-TextValue $ws "D5" '226.81'
Set-TextValue $ws "E5" '  -3.43%  '

# Row 6
Set-TextValue $ws "D6" '0.5795'
Set-TextValue $ws "E6" '  -3.10%  '

# Row 7
Set-TextValue $ws "D7" '1.002'
Set-TextValue $ws "E7" '  -0.12%  '

# Row 8
Set-TextValue $ws "D8" '0.2714'
Set-TextValue $ws "E8" '  -1.06%  '

# Row 9
Set-TextValue $ws "D9" '23.13'
Set-TextValue $ws "E9" '  -0.78%  '

# Row 10
Set-TextValue $ws "D10" '0.06601'

# Row 11
Set-TextValue $ws "D11" '0.07513'
Set-TextValue $ws "E11" '  -1.19%  '

# Row 12
Set-TextValue $ws "D12" '1.735.28'
Set-TextValue $ws "E12" '  -5.09%  '

# Row 13
Set-TextValue $ws "D13" '4.728'
Set-TextValue $ws "E13" '  -0.08%  '

# Row 14
Set-TextValue $ws "D14" '0.6051'
Set-TextValue $ws "E14" '  -3.08%  '

# Row 15
Set-TextValue $ws "D15" '1.983.01'
Set-TextValue $ws "E15" '  -4.41%  '

# Row 16
Set-TextValue $ws "D16" '74.23'
Set-TextValue $ws "E16" '  -3.74%  '

# Row 17
Set-TextValue $ws "D17" '0.000008661'
Set-TextValue $ws "E17" '  -11.49%  '

# Row 18
Set-TextValue $ws "D18" '28.040.58'
Set-TextValue $ws "E18" '  -2.91%  '

# Row 19
Set-TextValue $ws "D19" '5.330'
Set-TextValue $ws "E19" '  -3.97%  '

# Row 20
Set-TextValue $ws "D20" '1.002'
Set-TextValue $ws "E20" '  -0.17%  '

# Row 21
Set-TextValue $ws "D21" '205.12'
Set-TextValue $ws "E21" '  -4.61%  '

# Row 22
Set-TextValue $ws "E22" '  -1.92%  '

# Row 23
Set-TextValue $ws "D23" '6.630'
Set-TextValue $ws "E23" '  -3.49%  '

# Row 24
Set-TextValue $ws "D24" '1.003'
Set-TextValue $ws "E24" '  -0.15%  '

# Row 25
Set-TextValue $ws "D25" '149.87'
Set-TextValue $ws "E25" '  -4.15%  '

# Row 26
Set-TextValue $ws "D26" '8.036'
Set-TextValue $ws "E26" '  +1.46%  '

# Row 27
Set-TextValue $ws "D27" '0.1233'
Set-TextValue $ws "E27" '  -3.90%  '

# Row 28
Set-TextValue $ws "D28" '16.15'
Set-TextValue $ws "E28" '  -1.76%  '

# Row 29
Set-TextValue $ws "D29" '1.385'
Set-TextValue $ws "E29" '  -2.11%  '

# Row 30
Set-TextValue $ws "D30" '0.06184'
Set-TextValue $ws "E30" '  -5.15%  '

# Row 31
Set-TextValue $ws "E31" '  -3.39%  '

# Row 32
Set-TextValue $ws "D32" '3.735'
Set-TextValue $ws "E32" '  -2.48%  '

# Row 33
Set-TextValue $ws "D33" '3.716'
Set-TextValue $ws "E33" '  -1.31%  '

# Row 34
Set-TextValue $ws "D34" '1.678'
Set-TextValue $ws "E34" '  -2.65%  '

# Row 35
Set-TextValue $ws "D35" '1.036'
Set-TextValue $ws "E35" '  -5.09%  '

# Row 36
Set-TextValue $ws "D36" '0.6377'
Set-TextValue $ws "E36" '  -1.16%  '

# Row 37
Set-TextValue $ws "D37" '2.463'
Set-TextValue $ws "E37" '  -2.90%  '

# Row 38
Set-TextValue $ws "D38" '2.718'
Set-TextValue $ws "E38" '  -0.94%  '

# Row 39
Set-TextValue $ws "D39" '0.01675'
Set-TextValue $ws "E39" '  -4.54%  '

# Row 40
Set-TextValue $ws "D40" '1.127.25'
Set-TextValue $ws "E40" '  -0.96%  '

# Row 41
Set-TextValue $ws "D41" '6.218'
Set-TextValue $ws "E41" '  -4.21%  '

# Row 42
Set-TextValue $ws "D42" '0.8732'

# Row 43
Set-TextValue $ws "D43" '1.004'
Set-TextValue $ws "E43" '  -0.04%  '

# Row 44
Set-TextValue $ws "D44" '99.61'
Set-TextValue $ws "E44" '  -0.37%  '

# Row 45
Set-TextValue $ws "D45" '1.895.18'
Set-TextValue $ws "E45" '  -4.64%  '

# Row 46
Set-TextValue $ws "D46" '59.46'
Set-TextValue $ws "E46" '  -3.53%  '

# Row 47
Set-TextValue $ws "B47" 'BabyDogeCoin'
Set-TextValue $ws "C47" 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue $ws "D47" '0.00000000109'
Set-TextValue $ws "E47" '  -3.66%  '

# Row 48
Set-TextValue $ws "B48" 'RenderToken'
Set-TextValue $ws "C48" 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws "D48" '1.582'
Set-TextValue $ws "E48" '  -0.94%  '

# Row 49
Set-TextValue $ws "D49" '8.225'
Set-TextValue $ws "E49" '  -3.03%  '

# Row 50
Set-TextValue $ws "D50" '0.05380'
Set-TextValue $ws "E50" '  -2.31%  '

# Row 51
Set-TextValue $ws "D51" '6.284'
Set-TextValue $ws "E51" '  -1.81%  '
